{"js": "// Add a \"Citations\" section (heading text + a Neo4j reference line) right\n// before the final (bookmark-only) paragraph at the end of the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n// The very last paragraph in the document only carries the `_GoBack`\n// bookmark (no visible text) -- the new content is inserted right before it,\n// so the bookmark stays on the last paragraph of the document.\nconst lastParagraph = items[items.length - 1];\n\n// Each insertParagraph call with \"Before\" puts the new paragraph\n// immediately above the anchor paragraph. Insert \"Citations\" first so it\n// ends up above the bookmark paragraph, then insert the reference line\n// \"Before\" the (still last) bookmark paragraph so it lands between the two,\n// giving the final order:\n//   ... <empty paragraph> / Citations / \\[3\\] [Neo4j](https://neo4j.com/) / <bookmark paragraph>\nlastParagraph.insertParagraph(\"Citations\", Word.InsertLocation.before);\nlastParagraph.insertParagraph(\"\\\\[3\\\\] [Neo4j](https://neo4j.com/)\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Add a \"Citations\" section (heading text + a Neo4j reference line) right\n# before the final (bookmark-only) paragraph at the end of the document.\n$d = $word.ActiveDocument\n\n# The very last paragraph in the document only carries the `_GoBack`\n# bookmark (no visible text); insert the new paragraphs directly in front of\n# it so the bookmark stays on the document's last paragraph.\n$lastParagraph = $d.Paragraphs.Last\n$targetRange = $lastParagraph.Range\n\n# InsertBefore with embedded carriage returns creates one new paragraph per\n# line, landing right above the bookmark paragraph in this order:\n#   ... <empty paragraph> / Citations / \\[3\\] [Neo4j](https://neo4j.com/) / <bookmark paragraph>\n$targetRange.InsertBefore(\"Citations`r\" + \"\\[3\\] [Neo4j](https://neo4j.com/)`r\")\n"}
